$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 7: update Fecha, Volumen, Precio mínimo/máximo/promedio, Precio $/Kg
# ---------------------------------------------------------------------------
$ws.Cells.Item(7, 4).Value = 44421   # D7 Fecha
$ws.Cells.Item(7, 10).Value = 80     # J7 Volumen
$ws.Cells.Item(7, 11).Value = 16500  # K7 Precio minimo
$ws.Cells.Item(7, 12).Value = 16500  # L7 Precio maximo
$ws.Cells.Item(7, 13).Value = 16500  # M7 Precio promedio ponderado
$ws.Cells.Item(7, 16).Value = 550    # P7 Precio $/Kg

# ---------------------------------------------------------------------------
# Row 8: update Fecha, Volumen, Precio mínimo/máximo/promedio, Precio $/Kg
# ---------------------------------------------------------------------------
$ws.Cells.Item(8, 4).Value = 44400   # D8 Fecha
$ws.Cells.Item(8, 10).Value = 70     # J8 Volumen
$ws.Cells.Item(8, 11).Value = 15000  # K8 Precio minimo
$ws.Cells.Item(8, 12).Value = 15000  # L8 Precio maximo
$ws.Cells.Item(8, 13).Value = 15000  # M8 Precio promedio ponderado
$ws.Cells.Item(8, 16).Value = 500    # P8 Precio $/Kg

# ---------------------------------------------------------------------------
# Row 9: replaced entirely with new "Espanola" record (formerly held the
# "Madrigal" record, which now moves down to row 10)
# ---------------------------------------------------------------------------
$ws.Cells.Item(9, 1).Value = 4
$ws.Cells.Item(9, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(9, 3).Value = "Los Lagos"
$ws.Cells.Item(9, 4).Value = 44390
$ws.Cells.Item(9, 4).NumberFormat = $ws.Cells.Item(8, 4).NumberFormat
$ws.Cells.Item(9, 5).Value = 10
$ws.Cells.Item(9, 6).Value = 100112013
$ws.Cells.Item(9, 7).Value = "Alcachofa"
$ws.Cells.Item(9, 8).Value = "Española"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 80
$ws.Cells.Item(9, 11).Value = 16000
$ws.Cells.Item(9, 12).Value = 16000
$ws.Cells.Item(9, 13).Value = 16000
$ws.Cells.Item(9, 14).Value = "`$/caja 30 unidades"
$ws.Cells.Item(9, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(9, 16).Value = 533
$ws.Cells.Item(9, 17).Value = 30
$ws.Cells.Item(9, 18).Value = "Hortaliza"

# ---------------------------------------------------------------------------
# Row 10 (new): holds the "Madrigal" record that used to be in row 9
# ---------------------------------------------------------------------------
$ws.Cells.Item(10, 1).Value = 4
$ws.Cells.Item(10, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(10, 3).Value = "Los Lagos"
$ws.Cells.Item(10, 4).Value = 44161
$ws.Cells.Item(10, 4).NumberFormat = $ws.Cells.Item(8, 4).NumberFormat
$ws.Cells.Item(10, 5).Value = 10
$ws.Cells.Item(10, 6).Value = 100112013
$ws.Cells.Item(10, 7).Value = "Alcachofa"
$ws.Cells.Item(10, 8).Value = "Madrigal"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 30
$ws.Cells.Item(10, 11).Value = 11000
$ws.Cells.Item(10, 12).Value = 11000
$ws.Cells.Item(10, 13).Value = 11000
$ws.Cells.Item(10, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(10, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(10, 16).Value = 275
$ws.Cells.Item(10, 17).Value = 40
$ws.Cells.Item(10, 18).Value = "Hortaliza"

# ---------------------------------------------------------------------------
# Row 11 (new): brand new "Espanola" record
# ---------------------------------------------------------------------------
$ws.Cells.Item(11, 1).Value = 4
$ws.Cells.Item(11, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(11, 3).Value = "Los Lagos"
$ws.Cells.Item(11, 4).Value = 44418
$ws.Cells.Item(11, 4).NumberFormat = $ws.Cells.Item(8, 4).NumberFormat
$ws.Cells.Item(11, 5).Value = 10
$ws.Cells.Item(11, 6).Value = 100112013
$ws.Cells.Item(11, 7).Value = "Alcachofa"
$ws.Cells.Item(11, 8).Value = "Española"
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 80
$ws.Cells.Item(11, 11).Value = 16000
$ws.Cells.Item(11, 12).Value = 16000
$ws.Cells.Item(11, 13).Value = 16000
$ws.Cells.Item(11, 14).Value = "`$/caja 30 unidades"
$ws.Cells.Item(11, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(11, 16).Value = 533
$ws.Cells.Item(11, 17).Value = 30
$ws.Cells.Item(11, 18).Value = "Hortaliza"
